$wb = $excel.ActiveWorkbook

# Sheet "Application" (sheet2.xml): append 3 new IDs after the last row (A48 -> A49:A51)
$wsApp = $wb.Worksheets.Item("Application")
$appValues = @("A227305", "A227306", "A227307")
$appStartRow = 49
for ($i = 0; $i -lt $appValues.Length; $i++) {
    $wsApp.Cells.Item($appStartRow + $i, 1).Value = $appValues[$i]
}

# Sheet "Project" (sheet3.xml): append 1 new ID after the last row (A26 -> A27)
$wsProj = $wb.Worksheets.Item("Project")
$wsProj.Cells.Item(27, 1).Value = "P-860402"

# Sheet "Request" (sheet4.xml): append 3 new IDs after the last row (A35 -> A36:A38)
$wsReq = $wb.Worksheets.Item("Request")
$reqValues = @("R02176", "R02178", "R02180")
$reqStartRow = 36
for ($i = 0; $i -lt $reqValues.Length; $i++) {
    $wsReq.Cells.Item($reqStartRow + $i, 1).Value = $reqValues[$i]
}
